$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove all existing hyperlinks first; they will be re-added after data is rewritten.
$ws.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = '2025-10-22 12:37:40'
$ws.Range("B2").Value = '【ECシステム開発】販売データ分析・AI提案・競合監視を統合した販売支援システム構築'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5418284'
$ws.Range("G2").Value = 410
$ws.Range("H2").Value = '🔥AI,Ai ◆開発,システム開発'

# Row 3
$ws.Range("A3").Value = '2025-10-22 12:37:40'
$ws.Range("B3").Value = '【AIで開発生産性を革新】AI活用推進エンジニア募集(副業・業務委託)'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5417964'
$ws.Range("G3").Value = 375
$ws.Range("H3").Value = '🔥AI,Ai ◆開発'

# Row 4
$ws.Range("A4").Value = '2025-10-22 12:37:40'
$ws.Range("B4").Value = '【業務委託/副業可】AI SaaS開発を牽引するCTO候補を募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5417967'
$ws.Range("G4").Value = 375
$ws.Range("H4").Value = '🔥AI,Ai ◆開発'

# Row 5
$ws.Range("A5").Value = '2025-10-22 12:37:40'
$ws.Range("B5").Value = '【継続依頼あり】AI×業務効率化のスペシャリスト募集!'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5418075'
$ws.Range("G5").Value = 373
$ws.Range("H5").Value = '🔥AI,Ai ◆効率化'

# Row 6
$ws.Range("A6").Value = '2025-10-22 12:37:40'
$ws.Range("B6").Value = '【AI技術顧問/戦略アドバイザー募集】最先端AIで事業の非連続な成長を牽引するエキスパート求む'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5417960'
$ws.Range("G6").Value = 310
$ws.Range("H6").Value = '🔥AI,Ai'

# Row 7
$ws.Range("A7").Value = '2025-10-22 12:37:40'
$ws.Range("B7").Value = 'GASと生成AIを活用したスプレッドシートの作り方レクチャー'
$ws.Range("C7").Value = 'システム開発'
$ws.Range("D7").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E7").Value = '期限情報なし'
$ws.Range("F7").Value = 'https://www.lancers.jp/work/detail/5418291'
$ws.Range("G7").Value = 298
$ws.Range("H7").Value = '🔥AI,Ai'

# Row 8
$ws.Range("A8").Value = '2025-10-22 12:37:40'
$ws.Range("B8").Value = '【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪'
$ws.Range("C8").Value = 'システム開発'
$ws.Range("D8").Value = '20,000 円 ~ 50,000 円 / 固定'
$ws.Range("E8").Value = '期限情報なし'
$ws.Range("F8").Value = 'https://www.lancers.jp/work/detail/5217096'
$ws.Range("G8").Value = 243
$ws.Range("H8").Value = '🔥API ◆ツール'

# Row 9
$ws.Range("A9").Value = '2025-10-22 12:37:40'
$ws.Range("B9").Value = '【急募】キントーン見積書をエクセルに変換するツール開発'
$ws.Range("C9").Value = 'システム開発'
$ws.Range("D9").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E9").Value = '期限情報なし'
$ws.Range("F9").Value = 'https://www.lancers.jp/work/detail/5418067'
$ws.Range("G9").Value = 120
$ws.Range("H9").Value = '◆ツール,開発'

# Row 10
$ws.Range("A10").Value = '2025-10-22 12:37:40'
$ws.Range("B10").Value = '【急募】Webアプリ開発エンジニア募集!フルリモート可'
$ws.Range("C10").Value = 'システム開発'
$ws.Range("D10").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E10").Value = '期限情報なし'
$ws.Range("F10").Value = 'https://www.lancers.jp/work/detail/5411585'
$ws.Range("G10").Value = 93
$ws.Range("H10").Value = '◆開発 ◇アプリ'

# Row 11
$ws.Range("A11").Value = '2025-10-22 12:37:40'
$ws.Range("B11").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件_ワーカー'
$ws.Range("C11").Value = 'システム開発'
$ws.Range("D11").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E11").Value = '期限情報なし'
$ws.Range("F11").Value = 'https://www.lancers.jp/work/detail/5418318'
$ws.Range("G11").Value = 75
$ws.Range("H11").Value = '◆開発'

# Row 12
$ws.Range("A12").Value = '2025-10-22 12:37:40'
$ws.Range("B12").Value = '大手クレジットカード企業向け、Google Cloudを利用したアジャイル開発共通基盤案件'
$ws.Range("C12").Value = 'システム開発'
$ws.Range("D12").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E12").Value = '期限情報なし'
$ws.Range("F12").Value = 'https://www.lancers.jp/work/detail/5418320'
$ws.Range("G12").Value = 75
$ws.Range("H12").Value = '◆開発'

# Row 13
$ws.Range("A13").Value = '2025-10-22 12:37:40'
$ws.Range("B13").Value = 'IISIA公式サイト WordPressアップデート&AWS運用整備 実施要領書'
$ws.Range("C13").Value = 'システム開発'
$ws.Range("D13").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E13").Value = '期限情報なし'
$ws.Range("F13").Value = 'https://www.lancers.jp/work/detail/5418421'
$ws.Range("G13").Value = 65
$ws.Range("H13").Value = '◇サイト ○WordPress'

# Row 14
$ws.Range("A14").Value = '2025-10-22 12:37:40'
$ws.Range("B14").Value = '進行管理およびチームディレクションを担当'
$ws.Range("C14").Value = 'システム開発'
$ws.Range("D14").Value = '~ 5,000 円 / 固定'
$ws.Range("E14").Value = '期限情報なし'
$ws.Range("F14").Value = 'https://www.lancers.jp/work/detail/5418064'
$ws.Range("G14").Value = 30
$ws.Range("H14").Value = '◇管理'

# Row 15
$ws.Range("A15").Value = '2025-10-22 12:37:40'
$ws.Range("B15").Value = '自社HPに見積自動受付システムを設置したい'
$ws.Range("C15").Value = 'システム開発'
$ws.Range("D15").Value = '200,000 円 ~ 300,000 円 / 固定'
$ws.Range("E15").Value = '期限情報なし'
$ws.Range("F15").Value = 'https://www.lancers.jp/work/detail/5418456'
$ws.Range("G15").Value = 33

# Row 16
$ws.Range("A16").Value = '2025-10-22 12:37:40'
$ws.Range("B16").Value = '【WP安全アップデート+AWS運用まで一括/haradatakeo.com(~45万円以下)】'
$ws.Range("C16").Value = 'システム開発'
$ws.Range("D16").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E16").Value = '期限情報なし'
$ws.Range("F16").Value = 'https://www.lancers.jp/work/detail/5418426'
$ws.Range("G16").Value = 25

# Row 17
$ws.Range("A17").Value = '2025-10-22 12:37:40'
$ws.Range("B17").Value = 'サブスクペイからCSVデータをダウンロードし、データベース同期するプログラムの作成'
$ws.Range("C17").Value = 'システム開発'
$ws.Range("D17").Value = '50,000 円 ~ 100,000 円 / 固定'
$ws.Range("E17").Value = '期限情報なし'
$ws.Range("F17").Value = 'https://www.lancers.jp/work/detail/5418241'
$ws.Range("G17").Value = 18

# Row 18
$ws.Range("A18").Value = '2025-10-22 12:37:40'
$ws.Range("B18").Value = 'ハードウェアの設定設置と保守サポート依頼|東京周辺対応可能な方'
$ws.Range("C18").Value = 'システム開発'
$ws.Range("D18").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E18").Value = '期限情報なし'
$ws.Range("F18").Value = 'https://www.lancers.jp/work/detail/5418084'
$ws.Range("G18").Value = 18

# Row 19
$ws.Range("A19").Value = '2025-10-22 12:37:40'
$ws.Range("B19").Value = '【急募】Meta広告のコンバージョン計測設定をお手伝いください!'
$ws.Range("C19").Value = 'システム開発'
$ws.Range("D19").Value = '5,000 円 ~ 10,000 円 / 固定'
$ws.Range("E19").Value = '期限情報なし'
$ws.Range("F19").Value = 'https://www.lancers.jp/work/detail/5418533'
$ws.Range("G19").Value = 10

# Re-add hyperlinks for the URL column (F), in row order, so relationship ids are issued sequentially.
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5418284', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418284') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5417964', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5417964') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5417967', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5417967') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5418075', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418075') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5417960', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5417960') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), 'https://www.lancers.jp/work/detail/5418291', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418291') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F8"), 'https://www.lancers.jp/work/detail/5217096', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5217096') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F9"), 'https://www.lancers.jp/work/detail/5418067', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418067') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F10"), 'https://www.lancers.jp/work/detail/5411585', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5411585') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F11"), 'https://www.lancers.jp/work/detail/5418318', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418318') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F12"), 'https://www.lancers.jp/work/detail/5418320', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418320') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F13"), 'https://www.lancers.jp/work/detail/5418421', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418421') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F14"), 'https://www.lancers.jp/work/detail/5418064', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418064') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F15"), 'https://www.lancers.jp/work/detail/5418456', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418456') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F16"), 'https://www.lancers.jp/work/detail/5418426', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418426') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F17"), 'https://www.lancers.jp/work/detail/5418241', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418241') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F18"), 'https://www.lancers.jp/work/detail/5418084', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418084') | Out-Null
$ws.Hyperlinks.Add($ws.Range("F19"), 'https://www.lancers.jp/work/detail/5418533', [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, 'https://www.lancers.jp/work/detail/5418533') | Out-Null

# Column width adjustments (B: 49 -> 52, H: 13 -> 19).
# Excel COM ColumnWidth setter pads by 5/6 of a character unit before it is
# persisted to the OOXML <col width> attribute, so subtract that offset here
# to land exactly on the target stored widths of 52 and 19.
$ws.Columns.Item(2).ColumnWidth = 52 - (5/6)
$ws.Columns.Item(8).ColumnWidth = 19 - (5/6)

Write-Output "edit complete"
